$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.933.37"
$ws.Range("E2").Value = "  -2.52%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.658.30"
$ws.Range("E3").Value = "  +2.74%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "407.40"
$ws.Range("E5").Value = "  -1.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.46"
$ws.Range("E6").Value = "  +3.95%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.648.76"
$ws.Range("E7").Value = "  +2.64%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.621"
$ws.Range("E8").Value = "  -3.74%  "

$ws.Range("E9").Value = "  +0.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.728"
$ws.Range("E10").Value = "  -5.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.162"
$ws.Range("E11").Value = "  -6.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000320"
$ws.Range("E12").Value = "  +0.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.22"
$ws.Range("E13").Value = "  +0.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.94"
$ws.Range("E14").Value = "  +1.50%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.263.32"
$ws.Range("E15").Value = "  +3.33%  "

$ws.Range("E16").Value = "  -1.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.74"
$ws.Range("E17").Value = "  +12.54%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.02"
$ws.Range("E18").Value = "  -0.13%  "

$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.627.55"
$ws.Range("E19").Value = "  +1.11%  "

$ws.Range("E20").Value = "  -2.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "64.894.09"
$ws.Range("E21").Value = "  -2.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "421.75"
$ws.Range("E22").Value = "  -4.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.30"
$ws.Range("E23").Value = "  +18.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.99"
$ws.Range("E24").Value = "  -3.38%  "

$ws.Range("E25").Value = "  -3.58%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "35.80"
$ws.Range("E26").Value = "  +4.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.21"
$ws.Range("E27").Value = "  -2.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.40"
$ws.Range("E28").Value = "  -5.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.13"
$ws.Range("E29").Value = "  +5.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.72"
$ws.Range("E30").Value = "  +4.04%  "

$ws.Range("E31").Value = "  -1.67%  "

$ws.Range("E32").Value = "  +2.20%  "

$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "41.62"
$ws.Range("E33").Value = "  +6.74%  "

$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.96"
$ws.Range("E34").Value = "  -3.82%  "

$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.161"
$ws.Range("E35").Value = "  +1.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.98"
$ws.Range("E36").Value = "  -0.59%  "

$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0465"
$ws.Range("E38").Value = "  -4.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.96"
$ws.Range("E39").Value = "  +31.43%  "

$ws.Range("E40").Value = "  -4.25%  "

$ws.Range("E41").Value = "  -0.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0650"
$ws.Range("E42").Value = "  -7.91%  "

$ws.Range("E43").Value = "  +4.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.86"
$ws.Range("E44").Value = "  +27.83%  "

$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.41"
$ws.Range("E45").Value = "  +3.28%  "

$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.08"
$ws.Range("E46").Value = "  +7.18%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.13"
$ws.Range("E47").Value = "  +23.07%  "

$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "143.68"
$ws.Range("E48").Value = "  -1.89%  "

$ws.Range("E49").Value = "  -4.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.52"
$ws.Range("E50").Value = "  -6.60%  "

$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.291"
$ws.Range("E51").Value = "  -4.35%  "

